# Applies the "Another set of samples collected, along with their details" edit
# to Collection-details.xlsx: updates Focus for rows 23/24/29 from "Normal" to
# "Postural", and appends 20 new sample rows (36-55) describing Josh35-Josh54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Focus column (E) for 3 existing rows: Normal -> Postural ---
$ws.Range('E23').Value = 'Postural'
$ws.Range('E24').Value = 'Postural'
$ws.Range('E29').Value = 'Postural'

# --- Append new sample rows 36-55 ---
# Row 36
$ws.Range('A36').Value = 'Josh35'
$ws.Range('B36').Value = 'Sitting at computer, using computer, then standing up and walking away half way through'
$ws.Range('C36').Value = 'Sitting/Walking'
$ws.Range('D36').Value = 'Hands active (moving)'
$ws.Range('E36').Value = 'Postural'
$ws.Range('F36').Value = 'Sitting (hands at rest)'
$ws.Range('G36').Value = 'Walking (hands at rest)'

# Row 37
$ws.Range('A37').Value = 'Josh36'
$ws.Range('B37').Value = 'Walking around room, collecting things, then sitting down at computer halfway and using computer'
$ws.Range('C37').Value = 'Walking/Sitting'
$ws.Range('D37').Value = 'Hands active (moving)'
$ws.Range('E37').Value = 'Postural'
$ws.Range('F37').Value = 'Walking (hands active)'
$ws.Range('G37').Value = 'Sitting (hands active)'

# Row 38
$ws.Range('A38').Value = 'Josh37'
$ws.Range('B38').Value = 'Laying down, hands at rest on chest, then sitting up, hands at side at rest'
$ws.Range('C38').Value = 'Laying/Sitting'
$ws.Range('D38').Value = 'Hands at rest'
$ws.Range('E38').Value = 'Postural'
$ws.Range('F38').Value = 'Laying (hands at rest)'
$ws.Range('G38').Value = 'Sitting (hands at rest)'

# Row 39
$ws.Range('A39').Value = 'Josh38'
$ws.Range('B39').Value = 'Sitting down, using mobile phone (swiping), then standing up and using mobile phone (swiping)'
$ws.Range('C39').Value = 'Sitting/Standing'
$ws.Range('D39').Value = 'Hands active (moving)'
$ws.Range('E39').Value = 'Postural'
$ws.Range('F39').Value = 'Sitting (hands active)'
$ws.Range('G39').Value = 'Standing (hands active)'

# Row 40
$ws.Range('A40').Value = 'Josh39'
$ws.Range('B40').Value = 'Laying down, using mobile phone (swiping), then standing up and walking while looking at phone (holding )'
$ws.Range('C40').Value = 'Laying/Walking'
$ws.Range('D40').Value = 'Hands active (moving)'
$ws.Range('E40').Value = 'Postural'
$ws.Range('F40').Value = 'Laying (hands active)'
$ws.Range('G40').Value = 'Walking (hands active)'

# Row 41
$ws.Range('A41').Value = 'Josh40'
$ws.Range('B41').Value = 'Sitting down, using computer when suddenly falling out of chair'
$ws.Range('C41').Value = 'Sitting/Falling'
$ws.Range('D41').Value = 'Hands active (moving)'
$ws.Range('E41').Value = 'Postural'
$ws.Range('F41').Value = 'Sitting (hands active)'
$ws.Range('G41').Value = 'Laying (hands at rest)'

# Row 42
$ws.Range('A42').Value = 'Josh41'
$ws.Range('B42').Value = 'Standing still, hands at side (at rest) when suddenly faillng to ground'
$ws.Range('C42').Value = 'Standing/Falling'
$ws.Range('D42').Value = 'Hands at side/rest'
$ws.Range('E42').Value = 'Postural'
$ws.Range('F42').Value = 'Sitting (hands at rest)'
$ws.Range('G42').Value = 'Laying (hands at rest)'

# Row 43
$ws.Range('A43').Value = 'Josh42'
$ws.Range('B43').Value = 'Walking hands at side, tripping and falling to ground'
$ws.Range('C43').Value = 'Walking/Falling'
$ws.Range('D43').Value = 'Hands at side/rest'
$ws.Range('E43').Value = 'Postural'
$ws.Range('F43').Value = 'Walking (hands at rest)'
$ws.Range('G43').Value = 'Laying (hands at rest)'

# Row 44
$ws.Range('A44').Value = 'Josh43'
$ws.Range('B44').Value = 'Standing using mobile phone (swiping), then sitting down in chair using mobile phone (swiping) half way through'
$ws.Range('C44').Value = 'Standing/Sitting'
$ws.Range('D44').Value = 'Hands active (moving)'
$ws.Range('E44').Value = 'Postural'
$ws.Range('F44').Value = 'Standing (hands active)'
$ws.Range('G44').Value = 'Sitting (hands active)'

# Row 45
$ws.Range('A45').Value = 'Josh44'
$ws.Range('B45').Value = 'Walking using mobile phone (swiping), then Laying down using mobile phoe (swiping) half way through'
$ws.Range('C45').Value = 'Walking/Laying'
$ws.Range('D45').Value = 'Hands active (moving)'
$ws.Range('E45').Value = 'Postural'
$ws.Range('F45').Value = 'Walking (hands active)'
$ws.Range('G45').Value = 'Laying (hands active)'

# Row 46
$ws.Range('A46').Value = 'Josh45'
$ws.Range('B46').Value = 'Walking using mobile phone (swiping), then stopping and standing still, which using mobile phone (swiping)'
$ws.Range('C46').Value = 'Walking/Standing'
$ws.Range('D46').Value = 'Hands active (moving)'
$ws.Range('E46').Value = 'Postural'
$ws.Range('F46').Value = 'Walking (hands active)'
$ws.Range('G46').Value = 'Standing (hands active)'

# Row 47
$ws.Range('A47').Value = 'Josh46'
$ws.Range('B47').Value = 'Standing with hands at side (at rest), then sitting down in chair with hands on armsrests (at rest)'
$ws.Range('C47').Value = 'Standing/Sitting'
$ws.Range('D47').Value = 'Hands at rest'
$ws.Range('E47').Value = 'Postural'
$ws.Range('F47').Value = 'Stairs (hands at rest)'
$ws.Range('G47').Value = 'Sitting (hands at rest)'

# Row 48
$ws.Range('A48').Value = 'Josh47'
$ws.Range('B48').Value = 'Standing with hands at side (at rest), then laying down on bed with hands on chest (at rest)'
$ws.Range('C48').Value = 'Standing/Laying'
$ws.Range('D48').Value = 'Hands at rest'
$ws.Range('E48').Value = 'Postural'
$ws.Range('F48').Value = 'Standing (hands at rest)'
$ws.Range('G48').Value = 'Laying (hands at rest)'

# Row 49
$ws.Range('A49').Value = 'Josh48'
$ws.Range('B49').Value = 'Standing using mobile phone (swiping), then laying down on bed using mobile phone (swiping) half way through'
$ws.Range('C49').Value = 'Standing/Laying'
$ws.Range('D49').Value = 'Hands active (moving)'
$ws.Range('E49').Value = 'Postural'
$ws.Range('F49').Value = 'Standing (hands active)'
$ws.Range('G49').Value = 'Laying (hands active)'

# Row 50
$ws.Range('A50').Value = 'Josh49'
$ws.Range('B50').Value = 'Standing still talking, hands active (gesticulating) when suddenly falling to ground halfway through'
$ws.Range('C50').Value = 'Standing/Falling'
$ws.Range('D50').Value = 'Hands active (moving)'
$ws.Range('E50').Value = 'Postural'
$ws.Range('F50').Value = 'Standing (hands active)'
$ws.Range('G50').Value = 'Laying (hands at rest)'

# Row 51
$ws.Range('A51').Value = 'Josh50'
$ws.Range('B51').Value = 'Laying down hands at side (at rest), then standing up and walking while hands at side (at rest)'
$ws.Range('C51').Value = 'Laying/Walking'
$ws.Range('D51').Value = 'Hands at rest'
$ws.Range('E51').Value = 'Postural'
$ws.Range('F51').Value = 'Laying (hands at rest)'
$ws.Range('G51').Value = 'Walking (hands at rest)'

# Row 52
$ws.Range('A52').Value = 'Josh51'
$ws.Range('B52').Value = 'Walking around room hands at side (at rest), then sitting down at computer halfway with hands at rest'
$ws.Range('C52').Value = 'Walking/Sitting'
$ws.Range('D52').Value = 'Hands at rest'
$ws.Range('E52').Value = 'Postural'
$ws.Range('F52').Value = 'Walking (hands at rest)'
$ws.Range('G52').Value = 'Sitting (hands at rest)'

# Row 53
$ws.Range('A53').Value = 'Josh52'
$ws.Range('B53').Value = 'Sitting at computer, using computer, then a subtle tonic condition occurs halfway through'
$ws.Range('C53').Value = 'Sitting'
$ws.Range('D53').Value = 'Hands active (moving)'
$ws.Range('E53').Value = 'Tonic'
$ws.Range('F53').Value = 'Sitting (hands active)'
$ws.Range('G53').Value = 'Sitting (hands at rest)'

# Row 54
$ws.Range('A54').Value = 'Josh53'
$ws.Range('B54').Value = 'Sitting at computer, using computer, then a strong tonic condition occurs halfway through'
$ws.Range('C54').Value = 'Sitting'
$ws.Range('D54').Value = 'Hands active (moving)'
$ws.Range('E54').Value = 'Tonic'
$ws.Range('F54').Value = 'Sitting (hands active)'
$ws.Range('G54').Value = 'Sitting (hands at rest)'

# Row 55
$ws.Range('A55').Value = 'Josh54'
$ws.Range('B55').Value = 'Sitting at computer, using computer, then a extreme tonic condition occurs halfway through'
$ws.Range('C55').Value = 'Sitting'
$ws.Range('D55').Value = 'Hands active (moving)'
$ws.Range('E55').Value = 'Tonic'
$ws.Range('F55').Value = 'Sitting (hands active)'
$ws.Range('G55').Value = 'Sitting (hands at rest)'

# --- Match the column B width Excel computed after the new (longer) text was added ---
$ws.Columns.Item(2).ColumnWidth = 90.66666666666667

# --- Restore the view state: scrolled down a bit, with the last new row selected ---
$null = $ws.Range('C55:G55').Select()
